$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$fullText = $tr.Text

# Fix the capitalization of "Openvswitch" -> "OpenVswitch"
$oldWord = "Openvswitch"
$newWord = "OpenVswitch"
$idx = $fullText.IndexOf($oldWord)
if ($idx -ge 0) {
    $c = $tr.Characters($idx + 1, $oldWord.Length)
    $c.Text = $newWord
}

# Refresh text after the rename above, then split the " Agent" run into
# a " " run and an "Agent" run (same visible text, but two separate runs).
$fullText = $tr.Text
$target = " Agent"
$idx2 = $fullText.IndexOf($target)
if ($idx2 -ge 0) {
    $spaceRange = $tr.Characters($idx2 + 1, 1)
    $spaceRange.Text = " "
    $agentRange = $tr.Characters($idx2 + 2, 5)
    $agentRange.Text = "Agent"
}
